# Generate Report for Handback
#
# The handback status report recomputes a handful of timestamps (and a
# status code) when the CI job regenerates the localization handoff/handback
# report. In the source workbook several rows share identical placeholder
# values (this is generated test fixture data), so updating a value updates
# every row/sheet cell that held that same value.
#
# Changes:
#   "2016-08-26 10:13:28" -> "2016-08-26 10:14:20"
#   "ht"                  -> "mt"
#   "2016-08-26 10:13:23" -> "2016-08-26 10:14:16"
#   "2016-08-26 10:13:40" -> "2016-08-26 10:14:32"
#   "2016-08-26 10:13:47" -> "2016-08-26 10:14:39"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 4 & 5 ---
$wsOverview.Range("G4").Value = "2016-08-26 10:14:20"
$wsOverview.Range("G5").Value = "2016-08-26 10:14:20"

# --- zh-cn sheet: "Status" (E), "Correspond Handoff Datetime" (H),
#     "Correspond Handback DateTime" (K) columns, rows 4 & 5 ---
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-26 10:14:16"
$wsZhCn.Range("H5").Value = "2016-08-26 10:14:16"
$wsZhCn.Range("K4").Value = "2016-08-26 10:14:32"
$wsZhCn.Range("K5").Value = "2016-08-26 10:14:32"

# --- de-de sheet: "Status" (E), "Correspond Handoff Datetime" (H),
#     "Correspond Handback DateTime" (K) columns, rows 4 & 5 ---
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-26 10:14:20"
$wsDeDe.Range("H5").Value = "2016-08-26 10:14:20"
$wsDeDe.Range("K4").Value = "2016-08-26 10:14:39"
$wsDeDe.Range("K5").Value = "2016-08-26 10:14:39"
